$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text-valued cells (Coin name, Link, Volume label) ---
$textUpdates = @{
    'B9' = 'WazirX'
    'C9' = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
    'E9' = '8WazirXWRX'
    'B10' = 'MandalaExchangeToken'
    'C10' = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
    'E10' = '9MandalaExchangeTokenMDX'
    'B11' = 'LiechtensteinCryptoassetsExchange'
    'C11' = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
    'E11' = '10LiechtensteinCryptoassetsExchangeLCX'
    'B12' = 'BitrueCoin'
    'C12' = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
    'E12' = '11BitrueCoinBTR'
    'B13' = 'BitMartToken'
    'C13' = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
    'E13' = '12BitMartTokenBMX'
    'B15' = 'BitForexToken'
    'C15' = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
    'E15' = '14BitForexTokenBF'
    'B16' = 'CoinExToken'
    'C16' = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
    'E16' = '15CoinExTokenCET'
    'B17' = 'One'
    'C17' = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
    'E17' = '16OneONE'
}
foreach ($addr in $textUpdates.Keys) {
    $ws.Range($addr).Value = $textUpdates[$addr]
}

# --- Numeric-looking price cells that must remain stored as text ---
# (the source data keeps prices as inline strings, e.g. "245.35",
#  so force text format before the write, then restore the original style)
$priceUpdates = @{
    'D2' = '245.35'
    'D3' = '23.06'
    'D4' = '5.404'
    'D5' = '0.06047'
    'D7' = '0.8080'
    'D8' = '0.9317'
    'D9' = '0.1426'
    'D10' = '0.07460'
    'D11' = '0.03363'
    'D12' = '0.03068'
    'D13' = '0.09359'
    'D14' = '3.941'
    'D15' = '0.001599'
    'D16' = '0.04837'
    'D17' = '0.0005943'
    'D18' = '0.005478'
    'D20' = '0.0009846'
    'D22' = '3.648'
    'D23' = '6.443'
    'D40' = '0.03974'
    'D41' = '0.006423'
    'D44' = '0.006283'
    'D45' = '0.00005225'
    'D49' = '0.002179'
}
foreach ($addr in $priceUpdates.Keys) {
    $rng = $ws.Range($addr)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $priceUpdates[$addr]
    $rng.Style = $origStyle
}
